# Update Yojimbo Profits leve-flipping market data (current avg prices, leve prices, profits)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I32").Value = 933.6667
$ws.Range("J32").Value = 1132.6666
$ws.Range("K32").Value = 933.6667
$ws.Range("L32").Value = 1132.6666
$ws.Range("M32").Value = -607.6667
$ws.Range("N32").Value = -1784.6666

$ws.Range("H129").Value = 825.88
$ws.Range("I129").Value = 478.8
$ws.Range("J129").Value = 912.65
$ws.Range("K129").Value = 1436.4
$ws.Range("L129").Value = 2737.95
$ws.Range("M129").Value = 3563.6
$ws.Range("N129").Value = -12737.95

$ws.Range("H141").Value = 5658.1816
$ws.Range("I141").Value = 6471.1113
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 19413.3339
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = -14233.3339
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2288.2368
$ws.Range("I61").Value = 1513.5416
$ws.Range("J61").Value = 3616.2856
$ws.Range("K61").Value = 1513.5416
$ws.Range("L61").Value = 3616.2856
$ws.Range("M61").Value = -1301.5416
$ws.Range("N61").Value = -4040.2856

$ws.Range("H132").Value = 2686.513
$ws.Range("I132").Value = 2044.0714
$ws.Range("J132").Value = 4321.8184
$ws.Range("K132").Value = 6132.2142
$ws.Range("L132").Value = 12965.4552
$ws.Range("M132").Value = -3602.2142
$ws.Range("N132").Value = -18025.4552

$ws.Range("H136").Value = 2288.2368
$ws.Range("I136").Value = 1513.5416
$ws.Range("J136").Value = 3616.2856
$ws.Range("K136").Value = 4540.6248
$ws.Range("L136").Value = 10848.8568
$ws.Range("M136").Value = -1990.6248
$ws.Range("N136").Value = -15948.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 335.25
$ws.Range("I22").Value = 335.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 335.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -162.25

$ws.Range("H86").Value = 3394.319
$ws.Range("I86").Value = 3203.3684
$ws.Range("J86").Value = 4200.5557
$ws.Range("K86").Value = 3203.3684
$ws.Range("L86").Value = 4200.5557
$ws.Range("M86").Value = -2080.3684
$ws.Range("N86").Value = -6446.5557

$ws.Range("H89").Value = 3394.319
$ws.Range("I89").Value = 3203.3684
$ws.Range("J89").Value = 4200.5557
$ws.Range("K89").Value = 16016.842
$ws.Range("L89").Value = 21002.7785
$ws.Range("M89").Value = -10400.842
$ws.Range("N89").Value = -32234.7785

$ws.Range("H94").Value = 802.03705
$ws.Range("I94").Value = 758.4783
$ws.Range("J94").Value = 1052.5
$ws.Range("K94").Value = 758.4783
$ws.Range("L94").Value = 1052.5
$ws.Range("M94").Value = -307.4783
$ws.Range("N94").Value = -1954.5

$ws.Range("H139").Value = 30000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 30000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -40280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 36172.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 36172.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 36172.5
$ws.Range("N23").Value = -36652.5

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

$ws.Range("H27").Value = 36172.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 36172.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 36172.5
$ws.Range("N27").Value = -36556.5

$ws.Range("H32").Value = 2500
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2500
$ws.Range("N32").Value = -3132

$ws.Range("H33").Value = 1799.8334
$ws.Range("I33").Value = 1639.8
$ws.Range("J33").Value = 2600
$ws.Range("K33").Value = 1639.8
$ws.Range("L33").Value = 2600
$ws.Range("M33").Value = -1260.8
$ws.Range("N33").Value = -3358

$ws.Range("H122").Value = 3372.5
$ws.Range("I122").Value = 6380.2856
$ws.Range("J122").Value = 1968.8667
$ws.Range("K122").Value = 19140.8568
$ws.Range("L122").Value = 5906.6001
$ws.Range("M122").Value = -16690.8568
$ws.Range("N122").Value = -10806.6001

$ws.Range("H133").Value = 36885.2
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 36885.2
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 36885.2
$ws.Range("N133").Value = -41945.2

$ws.Range("H141").Value = 36500
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 42000
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 42000
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -52360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 507.30768
$ws.Range("I92").Value = 460
$ws.Range("J92").Value = 536.875
$ws.Range("K92").Value = 1380
$ws.Range("L92").Value = 1610.625
$ws.Range("M92").Value = -132
$ws.Range("N92").Value = -4106.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 28261
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 28261
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 28261
$ws.Range("N39").Value = -29325

$ws.Range("H70").Value = 4717.75
$ws.Range("I70").Value = 4631
$ws.Range("J70").Value = 4928.4287
$ws.Range("K70").Value = 4631
$ws.Range("L70").Value = 4928.4287
$ws.Range("M70").Value = -4361
$ws.Range("N70").Value = -5468.4287

$ws.Range("H73").Value = 4717.75
$ws.Range("I73").Value = 4631
$ws.Range("J73").Value = 4928.4287
$ws.Range("K73").Value = 4631
$ws.Range("L73").Value = 4928.4287
$ws.Range("M73").Value = -3695
$ws.Range("N73").Value = -6800.4287

$ws.Range("H126").Value = 1731.5294
$ws.Range("I126").Value = 1728
$ws.Range("J126").Value = 1740
$ws.Range("K126").Value = 5184
$ws.Range("L126").Value = 5220
$ws.Range("M126").Value = -2714
$ws.Range("N126").Value = -10160

$ws.Range("H132").Value = 2154.0815
$ws.Range("I132").Value = 1994.8529
$ws.Range("J132").Value = 2515
$ws.Range("K132").Value = 5984.5587
$ws.Range("L132").Value = 7545
$ws.Range("M132").Value = -3454.5587
$ws.Range("N132").Value = -12605

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 530.75
$ws.Range("I22").Value = 483.6154
$ws.Range("J22").Value = 571.6
$ws.Range("K22").Value = 483.6154
$ws.Range("L22").Value = 571.6
$ws.Range("M22").Value = -188.6154
$ws.Range("N22").Value = -1161.6

$ws.Range("H27").Value = 530.75
$ws.Range("I27").Value = 483.6154
$ws.Range("J27").Value = 571.6
$ws.Range("K27").Value = 483.6154
$ws.Range("L27").Value = 571.6
$ws.Range("M27").Value = -376.6154
$ws.Range("N27").Value = -785.6

$ws.Range("H139").Value = 44946
$ws.Range("I139").Value = 21000.334
$ws.Range("J139").Value = 53925.625
$ws.Range("K139").Value = 21000.334
$ws.Range("L139").Value = 53925.625
$ws.Range("M139").Value = -15860.334
$ws.Range("N139").Value = -64205.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 294.75
$ws.Range("I113").Value = 262
$ws.Range("J113").Value = 327.5
$ws.Range("K113").Value = 786
$ws.Range("L113").Value = 982.5
$ws.Range("M113").Value = 1384
$ws.Range("N113").Value = -5322.5

$ws.Range("H136").Value = 920.2857
$ws.Range("I136").Value = 609.6923
$ws.Range("J136").Value = 1425
$ws.Range("K136").Value = 1829.0769
$ws.Range("L136").Value = 4275
$ws.Range("M136").Value = 720.9231
$ws.Range("N136").Value = -9375

$ws.Range("H139").Value = 61522.832
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 61522.832
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 61522.832
$ws.Range("N139").Value = -71802.83199999999
